$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.970.52"
$ws.Range("E2").Value = "  -0.06%  "
$ws.Range("D3").Value = "2.258.75"
$ws.Range("E3").Value = "  -0.34%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.50"
$ws.Range("E5").Value = "  +0.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.33"
$ws.Range("E6").Value = "  +1.87%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.522"
$ws.Range("E7").Value = "  -1.38%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.491"
$ws.Range("E9").Value = "  +0.82%  "
$ws.Range("E10").Value = "  +3.38%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0786"
$ws.Range("E11").Value = "  -1.71%  "
$ws.Range("E12").Value = "  +0.41%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.78"
$ws.Range("E13").Value = "  +1.48%  "
$ws.Range("D14").Value = "2.609.68"
$ws.Range("E14").Value = "  -0.25%  "
$ws.Range("E15").Value = "  +1.50%  "
$ws.Range("D16").Value = "2.250.54"
$ws.Range("E16").Value = "  -1.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.788"
$ws.Range("E17").Value = "  +0.01%  "
$ws.Range("D18").Value = "41.870.80"
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("E19").Value = "  -2.91%  "
$ws.Range("E20").Value = "  -1.51%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.97"
$ws.Range("E21").Value = "  -0.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.05"
$ws.Range("E22").Value = "  +0.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.39"
$ws.Range("E23").Value = "  -2.72%  "
$ws.Range("E24").Value = "  -1.65%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.90"
$ws.Range("E26").Value = "  -1.82%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.43"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "36.71"
$ws.Range("E28").Value = "  +3.61%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.11"
$ws.Range("E29").Value = "  +0.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.44"
$ws.Range("E30").Value = "  -2.51%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "160.30"
$ws.Range("E31").Value = "  -0.34%  "
$ws.Range("E32").Value = "  +0.03%  "
$ws.Range("E33").Value = "  -2.75%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.20"
$ws.Range("E34").Value = "  +3.57%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0734"
$ws.Range("E35").Value = "  -2.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.25"
$ws.Range("E36").Value = "  +0.84%  "
$ws.Range("E37").Value = "  -0.15%  "
$ws.Range("E38").Value = "  -3.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.113"
$ws.Range("E39").Value = "  -1.95%  "
$ws.Range("E40").Value = "  -0.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.96"
$ws.Range("E41").Value = "  -3.47%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.27"
$ws.Range("E42").Value = "  +0.45%  "
$ws.Range("D43").Value = "1.956.51"
$ws.Range("E43").Value = "  -2.65%  "
$ws.Range("E44").Value = "  -0.92%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.59"
$ws.Range("E45").Value = "  -4.82%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.98"
$ws.Range("E46").Value = "  -2.76%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.86"
$ws.Range("E47").Value = "  -1.94%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "52.81"
$ws.Range("E48").Value = "  -1.33%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "71.92"
$ws.Range("E49").Value = "  -0.44%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "91.33"
$ws.Range("E50").Value = "  -0.42%  "
$ws.Range("E51").Value = "  -1.59%  "
